$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-5
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212)
$newDate = 45212

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
